$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 563.5
$ws.Range("J2").Value = 327
$ws.Range("L2").Value = 327
$ws.Range("N2").Value = -553
$ws.Range("H17").Value = 1098.125
$ws.Range("I17").Value = 1444.4445
$ws.Range("J17").Value = 997.5806
$ws.Range("K17").Value = 4333.333500000001
$ws.Range("L17").Value = 2992.7418
$ws.Range("M17").Value = -4165.333500000001
$ws.Range("N17").Value = -3328.7418
$ws.Range("H32").Value = 733
$ws.Range("I32").Value = 798.6667
$ws.Range("J32").Value = 667.3333
$ws.Range("K32").Value = 798.6667
$ws.Range("L32").Value = 667.3333
$ws.Range("M32").Value = -472.6667
$ws.Range("N32").Value = -1319.3333
$ws.Range("H51").Value = 5181.8184
$ws.Range("J51").Value = 5666.6665
$ws.Range("L51").Value = 5666.6665
$ws.Range("N51").Value = -6634.6665
$ws.Range("H125").Value = 3327
$ws.Range("I125").Value = 2248.3333
$ws.Range("K125").Value = 20234.9997
$ws.Range("M125").Value = -17774.9997
$ws.Range("H132").Value = 5538.0586
$ws.Range("I132").Value = 4949.6
$ws.Range("J132").Value = 9951.5
$ws.Range("K132").Value = 14848.8
$ws.Range("L132").Value = 29854.5
$ws.Range("M132").Value = -12318.8
$ws.Range("N132").Value = -34914.5
$ws.Range("H135").Value = 1073.4
$ws.Range("I135").Value = 857.75
$ws.Range("K135").Value = 7719.75
$ws.Range("M135").Value = -5184.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6127.4443
$ws.Range("I32").Value = 3366.3333
$ws.Range("J32").Value = 10466.333
$ws.Range("K32").Value = 3366.3333
$ws.Range("L32").Value = 10466.333
$ws.Range("M32").Value = -3079.3333
$ws.Range("N32").Value = -11040.333
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060
$ws.Range("H137").Value = 41780
$ws.Range("J137").Value = 41780
$ws.Range("L137").Value = 41780
$ws.Range("N137").Value = -51980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27
$ws.Range("H59").Value = 118846
$ws.Range("J59").Value = 118846
$ws.Range("L59").Value = 118846
$ws.Range("N59").Value = -120540
$ws.Range("H94").Value = 2766.75
$ws.Range("I94").Value = 2451.5
$ws.Range("J94").Value = 3712.5
$ws.Range("K94").Value = 2451.5
$ws.Range("L94").Value = 3712.5
$ws.Range("M94").Value = -2000.5
$ws.Range("N94").Value = -4614.5
$ws.Range("H96").Value = 9000
$ws.Range("I96").Value = 9000
$ws.Range("K96").Value = 9000
$ws.Range("M96").Value = -6254
$ws.Range("H107").Value = 1335.0869
$ws.Range("I107").Value = 1359.3334
$ws.Range("J107").Value = 1247.8
$ws.Range("K107").Value = 1359.3334
$ws.Range("L107").Value = 1247.8
$ws.Range("M107").Value = 560.6666
$ws.Range("N107").Value = -5087.8
$ws.Range("H137").Value = 35366.668
$ws.Range("J137").Value = 40550
$ws.Range("L137").Value = 40550
$ws.Range("N137").Value = -50750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14289966
$ws.Range("I99").Value = 33336320
$ws.Range("J99").Value = 5200
$ws.Range("K99").Value = 33336320
$ws.Range("L99").Value = 5200
$ws.Range("M99").Value = -33334822
$ws.Range("N99").Value = -8196
$ws.Range("H126").Value = 14289966
$ws.Range("I126").Value = 33336320
$ws.Range("J126").Value = 5200
$ws.Range("K126").Value = 100008960
$ws.Range("L126").Value = 15600
$ws.Range("M126").Value = -100006490
$ws.Range("N126").Value = -20540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1339.3077
$ws.Range("I5").Value = 426.25
$ws.Range("J5").Value = 2800.2
$ws.Range("K5").Value = 1278.75
$ws.Range("L5").Value = 8400.599999999999
$ws.Range("M5").Value = -1166.75
$ws.Range("N5").Value = -8624.599999999999
$ws.Range("H80").Value = 9799.200000000001
$ws.Range("J80").Value = 9799.200000000001
$ws.Range("L80").Value = 29397.6
$ws.Range("N80").Value = -31269.6
$ws.Range("H83").Value = 9799.200000000001
$ws.Range("J83").Value = 9799.200000000001
$ws.Range("L83").Value = 88192.8
$ws.Range("N83").Value = -97552.8
$ws.Range("H113").Value = 600.1070999999999
$ws.Range("I113").Value = 617.8
$ws.Range("J113").Value = 590.2778
$ws.Range("K113").Value = 1853.4
$ws.Range("L113").Value = 1770.8334
$ws.Range("M113").Value = 316.6000000000001
$ws.Range("N113").Value = -6110.8334
$ws.Range("H131").Value = 6173726.5
$ws.Range("J131").Value = 890.52563
$ws.Range("L131").Value = 2671.57689
$ws.Range("N131").Value = -12751.57689
$ws.Range("H132").Value = 1490.6129
$ws.Range("I132").Value = 768.55554
$ws.Range("J132").Value = 2490.3845
$ws.Range("K132").Value = 6916.99986
$ws.Range("L132").Value = 22413.4605
$ws.Range("M132").Value = -4386.99986
$ws.Range("N132").Value = -27473.4605
$ws.Range("H133").Value = 2848.3333
$ws.Range("I133").Value = 3900
$ws.Range("J133").Value = 2322.5
$ws.Range("K133").Value = 11700
$ws.Range("L133").Value = 6967.5
$ws.Range("M133").Value = -6640
$ws.Range("N133").Value = -17087.5
$ws.Range("H135").Value = 1339.3077
$ws.Range("I135").Value = 426.25
$ws.Range("J135").Value = 2800.2
$ws.Range("K135").Value = 3836.25
$ws.Range("L135").Value = 25201.8
$ws.Range("M135").Value = -1301.25
$ws.Range("N135").Value = -30271.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 34073.6
$ws.Range("J46").Value = 34073.6
$ws.Range("L46").Value = 34073.6
$ws.Range("N46").Value = -34385.6
$ws.Range("H107").Value = 6173490
$ws.Range("I107").Value = 545.8461
$ws.Range("K107").Value = 545.8461
$ws.Range("M107").Value = 1374.1539
$ws.Range("H122").Value = 1969.6757
$ws.Range("I122").Value = 1277.875
$ws.Range("J122").Value = 3246.8462
$ws.Range("K122").Value = 3833.625
$ws.Range("L122").Value = 9740.5386
$ws.Range("M122").Value = -1383.625
$ws.Range("N122").Value = -14640.5386
$ws.Range("H126").Value = 2925.09
$ws.Range("I126").Value = 2943.9695
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8831.908500000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -6361.908500000001
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2696.7646
$ws.Range("I132").Value = 1432.6522
$ws.Range("K132").Value = 4297.9566
$ws.Range("M132").Value = -1767.9566
$ws.Range("H137").Value = 45100
$ws.Range("J137").Value = 45100
$ws.Range("L137").Value = 45100
$ws.Range("N137").Value = -55300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1823.6364
$ws.Range("I16").Value = 1856
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1856
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1686
$ws.Range("N16").Value = -1840
$ws.Range("H122").Value = 3883.303
$ws.Range("I122").Value = 2513.45
$ws.Range("J122").Value = 5990.769
$ws.Range("K122").Value = 7540.349999999999
$ws.Range("L122").Value = 17972.307
$ws.Range("M122").Value = -5090.349999999999
$ws.Range("N122").Value = -22872.307
$ws.Range("H132").Value = 2953.4226
$ws.Range("I132").Value = 1253.1915
$ws.Range("J132").Value = 6283.0415
$ws.Range("K132").Value = 3759.5745
$ws.Range("L132").Value = 18849.1245
$ws.Range("M132").Value = -1229.5745
$ws.Range("N132").Value = -23909.1245

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 39750
$ws.Range("J92").Value = 39750
$ws.Range("L92").Value = 39750
$ws.Range("N92").Value = -44742
$ws.Range("H122").Value = 2951.9033
$ws.Range("I122").Value = 2000.75
$ws.Range("J122").Value = 3966.4666
$ws.Range("K122").Value = 6002.25
$ws.Range("L122").Value = 11899.3998
$ws.Range("M122").Value = -3552.25
$ws.Range("N122").Value = -16799.3998
$ws.Range("H136").Value = 1456.2106
$ws.Range("I136").Value = 598.1177
$ws.Range("J136").Value = 8750
$ws.Range("K136").Value = 1794.3531
$ws.Range("L136").Value = 26250
$ws.Range("M136").Value = 755.6469
$ws.Range("N136").Value = -31350
